$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-07-01 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-02 Wednesday", 2)

# Update table cell answers (using Table.Cell(row, col) to avoid ambiguity from duplicate text)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "31÷8=3, 7"
$t.Cell(1, 2).Range.Text = "37÷8=4, 5"
$t.Cell(1, 3).Range.Text = "99÷7=14, 1"
$t.Cell(1, 4).Range.Text = "93÷8=11, 5"
$t.Cell(1, 5).Range.Text = "36÷8=4, 4"

$t.Cell(5, 1).Range.Text = "49÷9=5, 4"
$t.Cell(5, 2).Range.Text = "68÷2=34, 0"
$t.Cell(5, 3).Range.Text = "86÷8=10, 6"
$t.Cell(5, 4).Range.Text = "31÷7=4, 3"
$t.Cell(5, 5).Range.Text = "83÷7=11, 6"

$t.Cell(9, 1).Range.Text = "24÷8=3, 0"
$t.Cell(9, 2).Range.Text = "61÷7=8, 5"
$t.Cell(9, 3).Range.Text = "98÷4=24, 2"
$t.Cell(9, 4).Range.Text = "80÷5=16, 0"
$t.Cell(9, 5).Range.Text = "15÷5=3, 0"

$t.Cell(13, 1).Range.Text = "13÷9=1, 4"
$t.Cell(13, 2).Range.Text = "76÷3=25, 1"
$t.Cell(13, 3).Range.Text = "89÷2=44, 1"
$t.Cell(13, 4).Range.Text = "24÷3=8, 0"
$t.Cell(13, 5).Range.Text = "76÷5=15, 1"

$t.Cell(17, 1).Range.Text = "77÷5=15, 2"
$t.Cell(17, 2).Range.Text = "47÷2=23, 1"
$t.Cell(17, 3).Range.Text = "77÷8=9, 5"
$t.Cell(17, 4).Range.Text = "19÷9=2, 1"
$t.Cell(17, 5).Range.Text = "34÷9=3, 7"
